$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.168506145477295
$ws.Range("B1").Value = 1.302209973335266
$ws.Range("C1").Value = 1.903980612754822
$ws.Range("D1").Value = 1.942407131195068
$ws.Range("E1").Value = 0.9256572127342224
